# Remove the trailing "Ver no Jupiter ..." / "(c) 2020 ..." footer
# paragraphs (plus the blank paragraph that separated them from the
# bibliography text above), as published by the site rebuild.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph by its stable text fragment.
$jupiterIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("Ver no Jupiter")) {
        $jupiterIndex = $i
        break
    }
}

if ($jupiterIndex -eq -1) {
    throw "Could not find the 'Ver no Jupiter' paragraph"
}

$pJupiter = $d.Paragraphs.Item($jupiterIndex)
$pCopyright = $d.Paragraphs.Item($jupiterIndex + 1)
$pBlank = $d.Paragraphs.Item($jupiterIndex - 1)

if (-not $pCopyright.Range.Text.Contains("luizeleno@usp.br")) {
    throw "Unexpected paragraph after 'Ver no Jupiter' paragraph"
}

# Delete the blank paragraph immediately before "Ver no Jupiter ...",
# the "Ver no Jupiter ..." paragraph itself, and the copyright/footer
# paragraph that follows it - all in one shot, including their
# paragraph marks.
$delRange = $d.Range($pBlank.Range.Start, $pCopyright.Range.End)
$delRange.Delete()

Write-Output "Deleted footer paragraphs; Paragraphs.Count is now $($d.Paragraphs.Count)"
